# GitHub action artefacts added
#
# Adds a new "Abstract Title" paragraph style (based on Normal, followed
# by Abstract) and tightens the space-before on the existing "Abstract"
# style from 15pt (300 twips) down to 5pt (100 twips).

$d = $word.ActiveDocument

# --- 1. New "AbstractTitle" style -----------------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
# wdColor is packed 0x00BBGGRR from hex RGB 345A8A
$abstractTitle.Font.Color = 0x34 + (0x5A * 256) + (0x8A * 65536)

# --- 2. Tighten spacing above the existing "Abstract" style ---------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

Write-Output "AbstractTitle style added; Abstract spacing-before updated"
